$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Type I error rate values for rows 2-7 (columns B:K)
$values = @{
    2 = @(0.236, 0.191, 0.267, 0.251, 0.113, 0.188, 0.166, 0.167, 0.22, 0.392)
    3 = @(0.297, 0.261, 0.345, 0.295, 0.13, 0.323, 0.275, 0.228, 0.345, 0.654)
    4 = @(0.428, 0.424, 0.464, 0.437, 0.158, 0.658, 0.571, 0.485, 0.6899999999999999, 0.923)
    5 = @(0.581, 0.581, 0.606, 0.587, 0.212, 0.862, 0.767, 0.667, 0.856, 0.983)
    6 = @(0.672, 0.675, 0.6870000000000001, 0.635, 0.247, 0.93, 0.886, 0.785, 0.9399999999999999, 0.995)
    7 = @(0.739, 0.791, 0.78, 0.751, 0.287, 0.968, 0.946, 0.875, 0.979, 0.999)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = 2 + $i   # Column B = 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$i]
    }
}

# Remove row 8 entirely (was nvec=50), shrinking the used range to A1:L7
$ws.Rows.Item(8).Delete()
